$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 34 ("@ format to Number:") to host the
# new "Percentage Text to Number:" example. This pushes everything from row 34
# downward by one row, matching the target layout.
$ws.Rows.Item(34).Insert()

# New label cell (shared-string backed, like its neighbours).
$ws.Cells.Item(34, 2).Value2 = "Percentage Text to Number:"

# New value cell: text "55.12%" parsed to the number 0.5512 with a percentage
# number format (builtin numFmtId 10 => "0.00%"), demonstrating that setting
# the datatype to Number now also understands percentage text.
$ws.Cells.Item(34, 3).NumberFormat = "0.00%"
$ws.Cells.Item(34, 3).Value2 = 0.5512

# Column B needs to be a bit wider to comfortably fit the new, longer label.
$ws.Columns.Item(2).ColumnWidth = 25.09
